# aggiornamento a 9/09 compreso
# Appends 8 new daily rows (2021-09-02 .. 2021-09-09) to the COVID tracking
# sheet, continuing directly after the existing last row (366 / 2021-09-01).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-column formatting (style/number-format) from the last
# existing data row (A366) down onto the new date cells so they keep the
# same "YYYY-MM-DD HH:MM:SS" look & border as the rest of column A.
$srcDateCell = $ws.Range("A366")
$newDateCells = $ws.Range("A367:A374")
$srcDateCell.Copy($newDateCells)

# date serials, nuovi pos. (B), somma mobile 7gg. (C), somma mobile 7gg. per 100mila abitanti (D)
$dates = @(44441, 44442, 44443, 44444, 44445, 44446, 44447, 44448)
$bvals = @(1, 1, 2, 0, 1, 2, 0, 0)
$cvals = @(10, 10, 7, 7, 5, 7, 7, 6)
$dvals = @(114.9954001839926, 114.9954001839926, 80.49678012879485, 80.49678012879485, 57.49770009199631, 80.49678012879485, 80.49678012879485, 68.99724011039559)

for ($i = 0; $i -lt 8; $i++) {
    $r = 367 + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = $bvals[$i]
    $ws.Cells.Item($r, 3).Value = $cvals[$i]
    $ws.Cells.Item($r, 4).Value = $dvals[$i]
}
